# Auto-generated Excel COM-interop edit script.
# Applies the odds updates for existing fixtures (rows 2-22) and appends
# the new row 23 fixture (Colombian Primera B: Boca Juniors de Cali vs Orsomarso).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated odds for existing rows (2-22) ---
$ws.Range("F2").Value = 2.82
$ws.Range("AA3").Value = 11
$ws.Range("AB3").Value = 970
$ws.Range("AC3").Value = 970
$ws.Range("AG3").Value = 970
$ws.Range("AH3").Value = 970
$ws.Range("F3").Value = 14.5
$ws.Range("G3").Value = 28
$ws.Range("I3").Value = 1.19
$ws.Range("K3").Value = 17.5
$ws.Range("N3").Value = 7.4
$ws.Range("V3").Value = 5.7
$ws.Range("X3").Value = 970
$ws.Range("AA4").Value = 44
$ws.Range("AE4").Value = 29
$ws.Range("F4").Value = 2.66
$ws.Range("G4").Value = 3.05
$ws.Range("I4").Value = 2.72
$ws.Range("K4").Value = 4
$ws.Range("O4").Value = 1.25
$ws.Range("Q4").Value = 1.73
$ws.Range("T4").Value = 1.62
$ws.Range("V4").Value = 1.58
$ws.Range("W4").Value = 1.49
$ws.Range("X4").Value = 970
$ws.Range("Y4").Value = 970
$ws.Range("L5").Value = 1.01
$ws.Range("R5").Value = 1.55
$ws.Range("T5").Value = 1.59
$ws.Range("AA6").Value = 26
$ws.Range("AB6").Value = 14.5
$ws.Range("AF6").Value = 26
$ws.Range("AJ6").Value = 70
$ws.Range("AN6").Value = 40
$ws.Range("AO6").Value = 15.5
$ws.Range("F6").Value = 3.75
$ws.Range("G6").Value = 3.8
$ws.Range("H6").Value = 2.16
$ws.Range("I6").Value = 2.18
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.84
$ws.Range("W6").Value = 1.35
$ws.Range("Z6").Value = 13
$ws.Range("AA7").Value = 42
$ws.Range("AF7").Value = 30
$ws.Range("AG7").Value = 18.5
$ws.Range("AH7").Value = 970
$ws.Range("AJ7").Value = 90
$ws.Range("AO7").Value = 34
$ws.Range("F7").Value = 3.1
$ws.Range("N7").Value = 2.94
$ws.Range("P7").Value = 1.66
$ws.Range("Q7").Value = 2.18
$ws.Range("S7").Value = 3.75
$ws.Range("T7").Value = 1.87
$ws.Range("Z7").Value = 17.5
$ws.Range("AH8").Value = 17.5
$ws.Range("AK8").Value = 27
$ws.Range("AL8").Value = 36
$ws.Range("O8").Value = 1.02
$ws.Range("R8").Value = 1.08
$ws.Range("S8").Value = 1.74
$ws.Range("G9").Value = 23
$ws.Range("I9").Value = 1000
$ws.Range("V9").Value = 1.4
$ws.Range("W9").Value = 1.22
$ws.Range("N10").Value = 2.66
$ws.Range("P10").Value = 1.64
$ws.Range("R10").Value = 1.19
$ws.Range("S10").Value = 2.52
$ws.Range("T10").Value = 1.05
$ws.Range("U10").Value = 1.05
$ws.Range("AC11").Value = 10.5
$ws.Range("G11").Value = 4
$ws.Range("J11").Value = 2.72
$ws.Range("P11").Value = 1.64
$ws.Range("Q11").Value = 2.3
$ws.Range("W11").Value = 1.33
$ws.Range("AD12").Value = 970
$ws.Range("F12").Value = 1.62
$ws.Range("U12").Value = 2.18
$ws.Range("X12").Value = 970
$ws.Range("Y12").Value = 970
$ws.Range("H13").Value = 2.38
$ws.Range("I13").Value = 2.58
$ws.Range("X13").Value = 19
$ws.Range("L14").Value = 1.31
$ws.Range("R14").Value = 1.4
$ws.Range("S14").Value = 3.05
$ws.Range("T14").Value = 1.76
$ws.Range("AH15").Value = 970
$ws.Range("U15").Value = 2.1
$ws.Range("P16").Value = 1.89
$ws.Range("Q16").Value = 1.9
$ws.Range("F17").Value = 1.74
$ws.Range("G17").Value = 2.22
$ws.Range("H17").Value = 3.15
$ws.Range("I17").Value = 5.8
$ws.Range("K17").Value = 9.2
$ws.Range("L17").Value = 1.01
$ws.Range("M17").Value = 1.03
$ws.Range("N17").Value = 2.66
$ws.Range("P17").Value = 1.77
$ws.Range("Q17").Value = 1.51
$ws.Range("R17").Value = 1.27
$ws.Range("S17").Value = 2.4
$ws.Range("T17").Value = 1.05
$ws.Range("U17").Value = 1.05
$ws.Range("V17").Value = 1.21
$ws.Range("W17").Value = 1.81
$ws.Range("AA18").Value = 75
$ws.Range("I18").Value = 3.9
$ws.Range("L18").Value = 1.34
$ws.Range("V18").Value = 1.34
$ws.Range("X18").Value = 17.5
$ws.Range("F19").Value = 5.9
$ws.Range("G19").Value = 6
$ws.Range("V19").Value = 2.58
$ws.Range("H20").Value = 1.98
$ws.Range("N20").Value = 5.5
$ws.Range("R20").Value = 1.59
$ws.Range("S20").Value = 2.6
$ws.Range("AH21").Value = 970
$ws.Range("G21").Value = 2.06
$ws.Range("H21").Value = 4.2
$ws.Range("M21").Value = 1.01
$ws.Range("N21").Value = 2.36
$ws.Range("R21").Value = 1.17
$ws.Range("W21").Value = 1.95
$ws.Range("AC22").Value = 9
$ws.Range("AJ22").Value = 120
$ws.Range("AL22").Value = 100
$ws.Range("F22").Value = 3.65
$ws.Range("G22").Value = 4.4
$ws.Range("H22").Value = 2.14
$ws.Range("I22").Value = 2.4
$ws.Range("K22").Value = 3.5
$ws.Range("N22").Value = 2.58
$ws.Range("P22").Value = 1.57
$ws.Range("Q22").Value = 2.42
$ws.Range("R22").Value = 1.17
$ws.Range("V22").Value = 1.73
$ws.Range("W22").Value = 1.3
$ws.Range("X22").Value = 11.5
$ws.Range("Y22").Value = 8.8

# --- New row 23: Colombian Primera B, Boca Juniors de Cali vs Orsomarso ---
$ws.Range("A23").Value = "Colombian Primera B"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "2026-02-17"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "22:10:00"
$ws.Range("D23").Value = "Boca Juniors de Cali"
$ws.Range("E23").Value = "Orsomarso"
$ws.Range("F23").Value = 2.24
$ws.Range("G23").Value = 2.38
$ws.Range("H23").Value = 3.65
$ws.Range("I23").Value = 4.2
$ws.Range("J23").Value = 3.05
$ws.Range("K23").Value = 3.4
$ws.Range("L23").Value = 1.01
$ws.Range("M23").Value = 1.11
$ws.Range("N23").Value = 2.56
$ws.Range("O23").Value = 1.54
$ws.Range("P23").Value = 1.5
$ws.Range("Q23").Value = 2.36
$ws.Range("R23").Value = 1.18
$ws.Range("S23").Value = 2.96
$ws.Range("T23").Value = 2.1
$ws.Range("U23").Value = 1.72
$ws.Range("V23").Value = 1.31
$ws.Range("W23").Value = 1.73
$ws.Range("X23").Value = 9.8
$ws.Range("Y23").Value = 12
$ws.Range("Z23").Value = 1000
$ws.Range("AA23").Value = 1000
$ws.Range("AB23").Value = 1000
$ws.Range("AC23").Value = 1000
$ws.Range("AD23").Value = 1000
$ws.Range("AE23").Value = 1000
$ws.Range("AF23").Value = 1000
$ws.Range("AG23").Value = 1000
$ws.Range("AH23").Value = 1000
$ws.Range("AI23").Value = 1000
$ws.Range("AJ23").Value = 1000
$ws.Range("AK23").Value = 1000
$ws.Range("AL23").Value = 1000
$ws.Range("AM23").Value = 1000
$ws.Range("AN23").Value = 1000
$ws.Range("AO23").Value = 1000
